# Applies the "Modif suivi & planning projet" edit:
#  - Reworks the first task table (rows 3-13): new task list, new "Réalisée"/
#    "Information potentiel" columns (F/G) with conditional formatting on F.
#  - Shifts the second table ("Test du système") down by 3 rows and expands
#    it with 3 new tasks + formatting (wrap text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Make room: push the "Test du système" block (row 11 onward) down
#    by 3 rows so it starts at row 14, matching the final layout.
# ------------------------------------------------------------------
$ws.Range("A11:A13").EntireRow.Insert()

# ------------------------------------------------------------------
# 2) First table header row (row 3): add F3/G3 headers, and turn on
#    wrap text for the "Nom de la tâche" header (D3) to match the new
#    wrapped body column below it.
# ------------------------------------------------------------------
$ws.Range("F3").Value2 = "Réalisée"
$ws.Range("G3").Value2 = "Information potentiel"
$ws.Range("F3").HorizontalAlignment = -4108
$ws.Range("F3").VerticalAlignment = -4108
$ws.Range("F3").Font.Bold = $true
$ws.Range("G3").HorizontalAlignment = -4108
$ws.Range("G3").VerticalAlignment = -4108
$ws.Range("G3").Font.Bold = $true
$ws.Range("D3").WrapText = $true
$ws.Range("D15").WrapText = $true

# ------------------------------------------------------------------
# 3) First table body (rows 4-13): new task names + recomputed "Tâche n°"
#    formulas following the IF(D<>"";ROW(C)-3;"-") pattern.
# ------------------------------------------------------------------
$tasks1 = @(
    "Mettre en place son environnement de travail",
    "Choisir un language de programmation",
    "Sélectionner un logiciel de développement adapté",
    "Réaliser le brouillon de l'algorithme",
    "Réaliser un diagramme de classes",
    "Réaliser un diagramme de séquence",
    "Mettre en place une base de donnée",
    "Mettre en place l'interface",
    "Coder le programme",
    "Optimiser le code"
)

for ($i = 0; $i -lt $tasks1.Length; $i++) {
    $r = 4 + $i
    $ws.Range("D$r").Value2 = $tasks1[$i]
    $ws.Range("D$r").WrapText = $true
    $ws.Range("D$r").HorizontalAlignment = -4108
    $ws.Range("D$r").VerticalAlignment = -4108
    $ws.Range("C$r").Formula = "=IF((D$r<>`"`"),ROW(C$r)-3,`"-`")"
    $ws.Range("E$r").Value2 = "-"
}

# "Réalisée" flag for the first task.
$ws.Range("F4").Value2 = "Oui"

# ------------------------------------------------------------------
# 4) Second table ("Test du système") now starts at row 14; its header
#    (row 15) formulas already point at C3/D3/E3 after the row insert.
#    Populate the (now longer) task list in rows 16-20 and extend the
#    "Tâche n°" formulas down to row 28.
# ------------------------------------------------------------------
$tasks2 = @(
    "Préparation de la fiche recette",
    "Vérification matériel",
    "Mise en place du matériel (ordinateur, afficheur, etc…)",
    "Test de l'application (fonctionnement bouton, bdd, envoi trame, etc…) + correction si nécessaire",
    "Communication entre les appareils"
)

for ($i = 0; $i -lt $tasks2.Length; $i++) {
    $r = 16 + $i
    $ws.Range("D$r").Value2 = $tasks2[$i]
    $ws.Range("D$r").WrapText = $true
    $ws.Range("D$r").HorizontalAlignment = -4108
    $ws.Range("D$r").VerticalAlignment = -4108
    $ws.Range("E$r").Value2 = "-"
}

# "Tâche n°" formulas for rows 16-20 (IF(...),ROW(A..),"-") pattern, as
# typed by hand in the source workbook.
$ws.Range("C16").Formula = '=IF((D16<>""),ROW(A1),"-")'
$ws.Range("C17").Formula = '=IF((D17<>""),ROW(A2),"-")'
$ws.Range("C18").Formula = '=IF((D18<>""),ROW(A3),"-")'
$ws.Range("C19").Formula = '=IF((D19<>""),ROW(A4),"-")'
$ws.Range("C20").Formula = '=IF((D20<>""),ROW(A5),"-")'

# Rows 21-28: continuation of the same formula pattern (kept empty of
# D/E content, same as the source workbook).
$ws.Range("C21").Formula = '=IF((D21<>""),ROW(A6),"-")'
$ws.Range("C22").Formula = '=IF((D22<>""),ROW(A8),"-")'
$ws.Range("C23").Formula = '=IF((D23<>""),ROW(A9),"-")'
$ws.Range("C24").Formula = '=IF((D24<>""),ROW(A10),"-")'
$ws.Range("C25").Formula = '=IF((D25<>""),ROW(A11),"-")'
$ws.Range("C26").Formula = '=IF((D26<>""),ROW(A12),"-")'
$ws.Range("C27").Formula = '=IF((D27<>""),ROW(A13),"-")'
$ws.Range("C28").Formula = '=IF((D28<>""),ROW(A14),"-")'

# Row 19 wraps to two lines in the source file (taller row).
$ws.Rows.Item(19).RowHeight = 30

# ------------------------------------------------------------------
# 5) Conditional formatting: highlight F column cells containing "Oui"
#    (Green, Accent 6, Lighter 60% fill - theme color 9 / tint ~0.6).
# ------------------------------------------------------------------
$rngF = $ws.Range("F1:F1048576")
$cond = $rngF.FormatConditions.Add(9, 0, "Oui")
$cond.Text = "Oui"
$cond.Formula1 = '=NOT(ISERROR(SEARCH("Oui",F1)))'
$cond.Interior.Color = 11854022

# ------------------------------------------------------------------
# 6) Column sizing to better fit the new content.
# ------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 8.75
$ws.Columns.Item(7).ColumnWidth = 19.65

# Final selection, matching the last-edited cell in the source workbook.
[void]$ws.Range("D16").Select()

$wb.Save()
